$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2723.2
$ws.Range("I11").Value = 2723.2
$ws.Range("K11").Value = 2723.2
$ws.Range("M11").Value = -2583.2
$ws.Range("H17").Value = 2518671.2
$ws.Range("J17").Value = 2567059
$ws.Range("L17").Value = 7701177
$ws.Range("N17").Value = -7701513
$ws.Range("H101").Value = 358.77777
$ws.Range("I101").Value = 479.33334
$ws.Range("J101").Value = 298.5
$ws.Range("K101").Value = 1438.00002
$ws.Range("L101").Value = 895.5
$ws.Range("M101").Value = 183.9999800000001
$ws.Range("N101").Value = -4139.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1823.7273
$ws.Range("J86").Value = 3152.1538
$ws.Range("L86").Value = 3152.1538
$ws.Range("N86").Value = -5398.1538
$ws.Range("H89").Value = 1823.7273
$ws.Range("J89").Value = 3152.1538
$ws.Range("L89").Value = 15760.769
$ws.Range("N89").Value = -26992.769
$ws.Range("H99").Value = 4624.75
$ws.Range("I99").Value = 4650.2
$ws.Range("J99").Value = 4497.5
$ws.Range("K99").Value = 4650.2
$ws.Range("L99").Value = 4497.5
$ws.Range("M99").Value = -3152.2
$ws.Range("N99").Value = -7493.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 614211.6
$ws.Range("I31").Value = 2501778
$ws.Range("J31").Value = 74906.92999999999
$ws.Range("K31").Value = 2501778
$ws.Range("L31").Value = 74906.92999999999
$ws.Range("M31").Value = -2501483
$ws.Range("N31").Value = -75496.92999999999
$ws.Range("H34").Value = 614211.6
$ws.Range("I34").Value = 2501778
$ws.Range("J34").Value = 74906.92999999999
$ws.Range("K34").Value = 2501778
$ws.Range("L34").Value = 74906.92999999999
$ws.Range("M34").Value = -2501576
$ws.Range("N34").Value = -75310.92999999999
$ws.Range("H58").Value = 2183.75
$ws.Range("I58").Value = 2122.6
$ws.Range("J58").Value = 2285.6667
$ws.Range("K58").Value = 2122.6
$ws.Range("L58").Value = 2285.6667
$ws.Range("M58").Value = -1919.6
$ws.Range("N58").Value = -2691.6667
$ws.Range("H86").Value = 2007295.2
$ws.Range("J86").Value = 11072
$ws.Range("L86").Value = 11072
$ws.Range("N86").Value = -13318
$ws.Range("H89").Value = 2007295.2
$ws.Range("J89").Value = 11072
$ws.Range("L89").Value = 55360
$ws.Range("N89").Value = -66592
$ws.Range("H136").Value = 2183.75
$ws.Range("I136").Value = 2122.6
$ws.Range("J136").Value = 2285.6667
$ws.Range("K136").Value = 6367.799999999999
$ws.Range("L136").Value = 6857.000100000001
$ws.Range("M136").Value = -3817.799999999999
$ws.Range("N136").Value = -11957.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 749.5
$ws.Range("I5").Value = 898
$ws.Range("J5").Value = 452.5
$ws.Range("K5").Value = 2694
$ws.Range("L5").Value = 1357.5
$ws.Range("M5").Value = -2582
$ws.Range("N5").Value = -1581.5
$ws.Range("H46").Value = 1703.4286
$ws.Range("I46").Value = 386
$ws.Range("J46").Value = 4997
$ws.Range("K46").Value = 1158
$ws.Range("L46").Value = 14991
$ws.Range("M46").Value = -1067
$ws.Range("N46").Value = -15173
$ws.Range("H74").Value = 5000
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 15000
$ws.Range("N74").Value = -17122
$ws.Range("H75").Value = 578.25
$ws.Range("J75").Value = 500
$ws.Range("L75").Value = 1500
$ws.Range("N75").Value = -3496
$ws.Range("H77").Value = 5000
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 45000
$ws.Range("N77").Value = -55608
$ws.Range("H78").Value = 578.25
$ws.Range("J78").Value = 500
$ws.Range("L78").Value = 4500
$ws.Range("N78").Value = -14484
$ws.Range("H135").Value = 749.5
$ws.Range("I135").Value = 898
$ws.Range("J135").Value = 452.5
$ws.Range("K135").Value = 8082
$ws.Range("L135").Value = 4072.5
$ws.Range("M135").Value = -5547
$ws.Range("N135").Value = -9142.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.25
$ws.Range("I2").Value = 65
$ws.Range("J2").Value = 31
$ws.Range("K2").Value = 65
$ws.Range("L2").Value = 31
$ws.Range("M2").Value = 48
$ws.Range("N2").Value = -257
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H114").Value = 57944.5
$ws.Range("J114").Value = 57944.5
$ws.Range("L114").Value = 57944.5
$ws.Range("N114").Value = -66622.5
$ws.Range("H132").Value = 22884.865
$ws.Range("I132").Value = 29129.719
$ws.Range("K132").Value = 87389.15700000001
$ws.Range("M132").Value = -84859.15700000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 9278.111000000001
$ws.Range("J97").Value = 9278.111000000001
$ws.Range("L97").Value = 9278.111000000001
$ws.Range("N97").Value = -11260.111
$ws.Range("H132").Value = 3144.4614
$ws.Range("I132").Value = 2616.4666
$ws.Range("K132").Value = 7849.399800000001
$ws.Range("M132").Value = -5319.399800000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 13004
$ws.Range("J15").Value = 13004
$ws.Range("L15").Value = 13004
$ws.Range("N15").Value = -13580
$ws.Range("H33").Value = 7999.5
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 7999.5
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 7999.5
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -8499.5
$ws.Range("H36").Value = 7999.5
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 7999.5
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 7999.5
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -8499.5
$ws.Range("H37").Value = 99943.5
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H126").Value = 3853.5
$ws.Range("I126").Value = 3853.5
$ws.Range("K126").Value = 11560.5
$ws.Range("M126").Value = -9090.5
$ws.Range("H132").Value = 2572.1177
$ws.Range("I132").Value = 2572.1177
$ws.Range("K132").Value = 7716.353099999999
$ws.Range("M132").Value = -5186.353099999999
$ws.Range("H136").Value = 557110.1
$ws.Range("I136").Value = 589828.4
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 1769485.2
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -1766935.2
$ws.Range("N136").Value = -7800
